$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.820.95'
$ws.Range('E2').Value = '  -0.70%  '

$ws.Range('D3').Value = '1.937.47'
$ws.Range('E3').Value = '  -0.96%  '

$ws.Range('D4').Formula = '''1.000'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Formula = '''243.56'
$ws.Range('E5').Value = '  -0.94%  '

$ws.Range('D6').Formula = '''1.001'
$ws.Range('E6').Value = '  +0.02%  '

$ws.Range('D7').Formula = '''0.4880'
$ws.Range('E7').Value = '  -0.47%  '

$ws.Range('D8').Formula = '''0.2945'
$ws.Range('E8').Value = '  -1.03%  '

$ws.Range('D9').Formula = '''0.06890'
$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('D10').Formula = '''19.29'
$ws.Range('E10').Value = '  +0.74%  '

$ws.Range('D11').Formula = '''104.77'
$ws.Range('E11').Value = '  -2.92%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.939.66'
$ws.Range('E12').Value = '  -0.83%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Formula = '''0.07791'
$ws.Range('E13').Value = '  +0.46%  '

$ws.Range('D14').Formula = '''5.347'
$ws.Range('E14').Value = '  -2.55%  '

$ws.Range('D15').Formula = '''0.7004'
$ws.Range('E15').Value = '  -1.33%  '

$ws.Range('D16').Formula = '''272.91'
$ws.Range('E16').Value = '  -3.59%  '

$ws.Range('D17').Value = '30.805.39'
$ws.Range('E17').Value = '  -0.86%  '

$ws.Range('D18').Formula = '''0.000007728'

$ws.Range('E19').Value = '  -1.71%  '

$ws.Range('B20').Value = 'BitDAO'
$ws.Range('C20').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D20').Formula = '''0.4890'
$ws.Range('E20').Value = '  -1.51%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Formula = '''5.610'
$ws.Range('E21').Value = '  +1.19%  '

$ws.Range('D22').Formula = '''1.000'
$ws.Range('E22').Value = '  +0.03%  '

$ws.Range('B23').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D23').Value = '2.184.45'
$ws.Range('E23').Value = '  -1.07%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Formula = '''1.000'
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Formula = '''6.530'
$ws.Range('E25').Value = '  +0.01%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Formula = '''9.831'
$ws.Range('E26').Value = '  -0.29%  '

$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Formula = '''165.82'
$ws.Range('E27').Value = '  -2.17%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Formula = '''19.63'
$ws.Range('E28').Value = '  -2.22%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Formula = '''2.158'
$ws.Range('E29').Value = '  -3.53%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Formula = '''0.1038'
$ws.Range('E30').Value = '  -1.91%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Formula = '''1.390'
$ws.Range('E31').Value = '  -2.64%  '

$ws.Range('D32').Formula = '''4.585'
$ws.Range('E32').Value = '  -0.11%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Formula = '''1.557'
$ws.Range('E33').Value = '  -2.16%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Formula = '''4.380'
$ws.Range('E34').Value = '  -2.68%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Formula = '''0.04885'
$ws.Range('E35').Value = '  -2.09%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Formula = '''0.7592'
$ws.Range('E36').Value = '  -0.24%  '

$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Formula = '''1.150'
$ws.Range('E37').Value = '  -2.99%  '

$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Formula = '''0.9998'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').Formula = '''2.731'
$ws.Range('E39').Value = '  +0.02%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Formula = '''0.02010'
$ws.Range('E40').Value = '  -1.07%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Formula = '''80.34'
$ws.Range('E41').Value = '  +7.85%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Formula = '''2.657'
$ws.Range('E42').Value = '  -1.91%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Formula = '''6.497'
$ws.Range('E43').Value = '  -1.08%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Formula = '''2.081'
$ws.Range('E44').Value = '  -4.09%  '

$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Formula = '''0.9040'
$ws.Range('E45').Value = '  +1.82%  '

$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Formula = '''0.4434'
$ws.Range('E46').Value = '  -1.97%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Formula = '''108.14'
$ws.Range('E47').Value = '  -1.31%  '

$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').Formula = '''1.001'
$ws.Range('E48').Value = '  +0.02%  '

$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').Formula = '''7.760'
$ws.Range('E49').Value = '  -4.65%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.000.41'
$ws.Range('E50').Value = '  +1.95%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Formula = '''0.1246'
$ws.Range('E51').Value = '  -2.02%  '

